# Updated cryptos list on Fri May 19 15:53:09 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row. Three pairs of rows also had their ranking shuffle, so their
# Coin name (column B) and Link (column C) are swapped to the coin that now
# occupies that row, in addition to a refreshed Price/Volume:
#   rows 14/15  - TRON <-> Litecoin
#   rows 41-43  - RenderToken -> MXToken -> RenderToken -> FraxShare cycle
#   rows 48/49  - PaxDollar <-> NEARProtocol
#
# Column D holds numeric-looking *text* (e.g. "1.000", "0.8710",
# "26.904.58" -- note some aren't even valid numbers). Assigning a plain
# numeric-looking string to .Value lets Excel silently coerce it into a
# Number (dropping trailing zeros / reformatting), so every such price is
# written with a leading apostrophe to force a literal text entry, exactly
# as typing '1.000 into Excel would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '26.904.58'
$ws.Cells.Item(2,5).Value = '  -1.02%  '

$ws.Cells.Item(3,4).Value = '1.816.68'
$ws.Cells.Item(3,5).Value = '  +0.15%  '

$ws.Cells.Item(4,4).Value = '''1.001'
$ws.Cells.Item(4,5).Value = '  -0.02%  '

$ws.Cells.Item(5,4).Value = '''310.02'
$ws.Cells.Item(5,5).Value = '  -0.68%  '

$ws.Cells.Item(6,4).Value = '''1.000'
$ws.Cells.Item(6,5).Value = '  -0.11%  '

$ws.Cells.Item(7,4).Value = '''0.4646'
$ws.Cells.Item(7,5).Value = '  +0.82%  '

$ws.Cells.Item(8,4).Value = '''0.3698'
$ws.Cells.Item(8,5).Value = '  -1.30%  '

$ws.Cells.Item(9,4).Value = '''0.07373'
$ws.Cells.Item(9,5).Value = '  -0.21%  '

$ws.Cells.Item(10,4).Value = '''0.8710'
$ws.Cells.Item(10,5).Value = '  +0.71%  '

$ws.Cells.Item(11,4).Value = '''20.44'
$ws.Cells.Item(11,5).Value = '  -0.72%  '

$ws.Cells.Item(12,4).Value = '1.820.55'
$ws.Cells.Item(12,5).Value = '  +0.38%  '

$ws.Cells.Item(13,4).Value = '''5.363'
$ws.Cells.Item(13,5).Value = '  -0.42%  '

# Row 14 was TRON, is now Litecoin (ranking shuffled)
$ws.Cells.Item(14,2).Value = 'Litecoin'
$ws.Cells.Item(14,3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(14,4).Value = '''92.19'
$ws.Cells.Item(14,5).Value = '  +0.44%  '

# Row 15 was Litecoin, is now TRON
$ws.Cells.Item(15,2).Value = 'TRON'
$ws.Cells.Item(15,3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(15,4).Value = '''0.07065'
$ws.Cells.Item(15,5).Value = '  -0.24%  '

$ws.Cells.Item(16,4).Value = '''6.510'
$ws.Cells.Item(16,5).Value = '  -2.20%  '

$ws.Cells.Item(17,5).Value = '  -0.08%  '

$ws.Cells.Item(18,4).Value = '''0.000008720'
$ws.Cells.Item(18,5).Value = '  -0.18%  '

$ws.Cells.Item(19,5).Value = '  +0.17%  '

$ws.Cells.Item(20,4).Value = '''14.71'
$ws.Cells.Item(20,5).Value = '  -1.09%  '

$ws.Cells.Item(21,4).Value = '26.966.27'
$ws.Cells.Item(21,5).Value = '  -0.81%  '

$ws.Cells.Item(22,4).Value = '''5.345'
$ws.Cells.Item(22,5).Value = '  +0.67%  '

$ws.Cells.Item(23,4).Value = '''10.54'
$ws.Cells.Item(23,5).Value = '  -3.33%  '

$ws.Cells.Item(24,4).Value = '2.099.32'
$ws.Cells.Item(24,5).Value = '  +2.92%  '

$ws.Cells.Item(25,4).Value = '''1.894'
$ws.Cells.Item(25,5).Value = '  -1.78%  '

$ws.Cells.Item(26,4).Value = '''151.78'
$ws.Cells.Item(26,5).Value = '  -0.05%  '

$ws.Cells.Item(27,4).Value = '''18.43'
$ws.Cells.Item(27,5).Value = '  -0.20%  '

$ws.Cells.Item(28,4).Value = '''2.175'
$ws.Cells.Item(28,5).Value = '  -2.09%  '

$ws.Cells.Item(29,4).Value = '''5.327'
$ws.Cells.Item(29,5).Value = '  +1.27%  '

$ws.Cells.Item(30,4).Value = '''115.53'
$ws.Cells.Item(30,5).Value = '  -1.20%  '

$ws.Cells.Item(31,4).Value = '''0.08902'
$ws.Cells.Item(31,5).Value = '  +0.16%  '

$ws.Cells.Item(32,4).Value = '''0.7620'
$ws.Cells.Item(32,5).Value = '  -1.25%  '

$ws.Cells.Item(33,4).Value = '''1.160'
$ws.Cells.Item(33,5).Value = '  -0.96%  '

$ws.Cells.Item(34,4).Value = '''4.489'
$ws.Cells.Item(34,5).Value = '  -0.54%  '

$ws.Cells.Item(35,4).Value = '''2.929'
$ws.Cells.Item(35,5).Value = '  +0.34%  '

$ws.Cells.Item(36,4).Value = '''0.9996'
$ws.Cells.Item(36,5).Value = '  -0.17%  '

$ws.Cells.Item(37,4).Value = '''1.098'
$ws.Cells.Item(37,5).Value = '  -2.13%  '

$ws.Cells.Item(38,4).Value = '''0.01959'
$ws.Cells.Item(38,5).Value = '  -0.05%  '

$ws.Cells.Item(39,4).Value = '''0.05258'
$ws.Cells.Item(39,5).Value = '  +0.41%  '

$ws.Cells.Item(40,4).Value = '''0.5374'
$ws.Cells.Item(40,5).Value = '  +1.77%  '

# Row 41 was RenderToken, is now MXToken (ranking shuffled)
$ws.Cells.Item(41,2).Value = 'MXToken'
$ws.Cells.Item(41,3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(41,4).Value = '''2.929'
$ws.Cells.Item(41,5).Value = '  +0.60%  '

# Row 42 was FraxShare, is now RenderToken
$ws.Cells.Item(42,2).Value = 'RenderToken'
$ws.Cells.Item(42,3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(42,4).Value = '''2.392'
$ws.Cells.Item(42,5).Value = '  +0.52%  '

# Row 43 was MXToken, is now FraxShare
$ws.Cells.Item(43,2).Value = 'FraxShare'
$ws.Cells.Item(43,3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(43,4).Value = '''7.231'
$ws.Cells.Item(43,5).Value = '  -0.41%  '

$ws.Cells.Item(44,4).Value = '''0.1660'
$ws.Cells.Item(44,5).Value = '  -1.16%  '

$ws.Cells.Item(45,4).Value = '''8.499'
$ws.Cells.Item(45,5).Value = '  -1.12%  '

$ws.Cells.Item(46,4).Value = '''0.4956'
$ws.Cells.Item(46,5).Value = '  -1.28%  '

$ws.Cells.Item(47,4).Value = '''10.41'
$ws.Cells.Item(47,5).Value = '  +0.63%  '

# Row 48 was PaxDollar, is now NEARProtocol (ranking shuffled)
$ws.Cells.Item(48,2).Value = 'NEARProtocol'
$ws.Cells.Item(48,3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(48,4).Value = '''1.676'
$ws.Cells.Item(48,5).Value = '  +0.42%  '

# Row 49 was NEARProtocol, is now PaxDollar
$ws.Cells.Item(49,2).Value = 'PaxDollar'
$ws.Cells.Item(49,3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(49,4).Value = '''0.9994'
$ws.Cells.Item(49,5).Value = '  -0.17%  '

$ws.Cells.Item(50,4).Value = '''103.08'
$ws.Cells.Item(50,5).Value = '  -1.94%  '

$ws.Cells.Item(51,4).Value = '''0.06291'
$ws.Cells.Item(51,5).Value = '  -0.41%  '
